$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so values like "1.000" or "19.11" are not
# auto-converted to numbers by Excel (matches original inlineStr/text storage).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2: Bitcoin
$ws.Range("D2").Value = '30.399.40'
$ws.Range("E2").Value = '  -0.01%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '1.927.44'
$ws.Range("E3").Value = '  +4.19%  '

# Row 4: TetherUSD
$ws.Range("D4").Value = '0.9993'
$ws.Range("E4").Value = '  -0.20%  '

# Row 5: BNB
$ws.Range("D5").Value = '240.73'
$ws.Range("E5").Value = '  +3.31%  '

# Row 6: USDC
$ws.Range("D6").Value = '0.9997'
$ws.Range("E6").Value = '  -0.14%  '

# Row 7: XRP
$ws.Range("D7").Value = '0.4768'
$ws.Range("E7").Value = '  +0.63%  '

# Row 8: Cardano
$ws.Range("D8").Value = '0.2871'
$ws.Range("E8").Value = '  +4.70%  '

# Row 9: Dogecoin
$ws.Range("D9").Value = '0.06598'
$ws.Range("E9").Value = '  +4.54%  '

# Row 10: Solana
$ws.Range("D10").Value = '19.11'
$ws.Range("E10").Value = '  +8.88%  '

# Row 11: Litecoin
$ws.Range("D11").Value = '107.15'
$ws.Range("E11").Value = '  +27.08%  '

# Row 12: WrappedEther
$ws.Range("D12").Value = '1.906.41'
$ws.Range("E12").Value = '  +3.14%  '

# Row 13: TRON
$ws.Range("D13").Value = '0.07618'
$ws.Range("E13").Value = '  +2.06%  '

# Row 14: Polkadot
$ws.Range("D14").Value = '5.129'
$ws.Range("E14").Value = '  +3.82%  '

# Row 15: Polygon
$ws.Range("D15").Value = '0.6576'
$ws.Range("E15").Value = '  +5.88%  '

# Row 16: BitcoinCash
$ws.Range("D16").Value = '302.37'
$ws.Range("E16").Value = '  +24.09%  '

# Row 17: WrappedBTC
$ws.Range("D17").Value = '30.417.99'
$ws.Range("E17").Value = '  +0.14%  '

# Row 18: Dai
$ws.Range("D18").Value = '1.000'
$ws.Range("E18").Value = '  -0.11%  '

# Row 19: Avalanche
$ws.Range("E19").Value = '  +2.59%  '

# Row 20: WrappedliquidstakedEther2.0
$ws.Range("D20").Value = '2.172.38'
$ws.Range("E20").Value = '  +3.96%  '

# Row 21: ShibaInu
$ws.Range("D21").Value = '0.000007494'
$ws.Range("E21").Value = '  +2.75%  '

# Row 22: Uniswap
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '5.294'
$ws.Range("E22").Value = '  +8.21%  '

# Row 23: BinanceUSD
$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D23").Value = '0.9993'
$ws.Range("E23").Value = '  -0.42%  '

# Row 24: Chainlink
$ws.Range("D24").Value = '6.302'
$ws.Range("E24").Value = '  +7.01%  '

# Row 25: Monero
$ws.Range("D25").Value = '168.08'
$ws.Range("E25").Value = '  +2.04%  '

# Row 26: Cosmos
$ws.Range("D26").Value = '9.213'
$ws.Range("E26").Value = '  +1.64%  '

# Row 27: EthereumClassic
$ws.Range("D27").Value = '20.03'
$ws.Range("E27").Value = '  +11.80%  '

# Row 28: LidoDAOToken
$ws.Range("D28").Value = '2.024'
$ws.Range("E28").Value = '  +8.81%  '

# Row 29: Stellar
$ws.Range("D29").Value = '0.1116'
$ws.Range("E29").Value = '  +8.30%  '

# Row 30: Toncoin
$ws.Range("D30").Value = '1.353'
$ws.Range("E30").Value = '  +0.52%  '

# Row 31: InternetComputer(DFINITY)
$ws.Range("D31").Value = '4.088'
$ws.Range("E31").Value = '  +1.60%  '

# Row 32: Filecoin
$ws.Range("E32").Value = '  +3.04%  '

# Row 33: Hedera
$ws.Range("D33").Value = '0.05001'
$ws.Range("E33").Value = '  +3.62%  '

# Row 34: ImmutableX
$ws.Range("D34").Value = '0.7410'
$ws.Range("E34").Value = '  +6.73%  '

# Row 35: ARBITRUM
$ws.Range("D35").Value = '1.151'
$ws.Range("E35").Value = '  +2.42%  '

# Row 36: HuobiToken
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = '2.734'
$ws.Range("E36").Value = '  +1.25%  '

# Row 37: VeChain
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '0.01941'
$ws.Range("E37").Value = '  +2.26%  '

# Row 38: MXToken
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").Value = '2.697'
$ws.Range("E38").Value = '  +0.72%  '

# Row 39: RenderToken
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '2.065'
$ws.Range("E39").Value = '  +3.46%  '

# Row 40: TrustWalletToken
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = '0.8766'
$ws.Range("E40").Value = '  +0.45%  '

# Row 41: Quant
$ws.Range("B41").Value = 'Quant'
$ws.Range("C41").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D41").Value = '107.10'
$ws.Range("E41").Value = '  +0.76%  '

# Row 42: FraxShare
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '5.820'
$ws.Range("E42").Value = '  +6.10%  '

# Row 43: Aave
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").Value = '69.98'
$ws.Range("E43").Value = '  +11.91%  '

# Row 44: PaxDollar
$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").Value = '0.9993'
$ws.Range("E44").Value = '  -0.16%  '

# Row 45: TheSandbox
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '0.4135'
$ws.Range("E45").Value = '  +2.36%  '

# Row 46: Aptos
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = '7.240'
$ws.Range("E46").Value = '  +1.56%  '

# Row 47: EnergySwap
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '9.301'
$ws.Range("E47").Value = '  +8.60%  '

# Row 48: Elrond
$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").Value = '34.86'
$ws.Range("E48").Value = '  +3.75%  '

# Row 49: Algorand
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = '0.1201'
$ws.Range("E49").Value = '  +0.57%  '

# Row 50: Cronos
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.05622'
$ws.Range("E50").Value = '  +1.98%  '

# Row 51: Decentraland
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").Value = '0.3841'
$ws.Range("E51").Value = '  +4.83%  '
